$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.424.76'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '1.673.67'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.81%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '221.09'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5362'
$ws.Range('E6').Value = '  +1.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.011'
$ws.Range('E7').Value = '  +0.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2675'
$ws.Range('E8').Value = '  +2.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06417'
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.06'
$ws.Range('E10').Value = '  +3.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07859'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.576'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '1.678.33'
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').Value = '1.902.85'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5664'
$ws.Range('E15').Value = '  +3.22%  '
$ws.Range('D16').Value = '0.0₅8210'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.39'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').Value = '26.467.56'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.716'
$ws.Range('E20').Value = '  +2.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '197.29'
$ws.Range('E21').Value = '  +3.30%  '
$ws.Range('E22').Value = '  +2.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.078'
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.45'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1234'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.276'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.23'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.514'
$ws.Range('E29').Value = '  +3.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05904'
$ws.Range('E30').Value = '  +2.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.292'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.589'
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.315'
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.627'
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9736'
$ws.Range('E35').Value = '  +2.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.854'
$ws.Range('E36').Value = '  +1.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.429'
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5840'
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01615'
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('D40').Value = '1.080.50'
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.917'
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8673'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '104.53'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').Value = '1.811.88'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.47'
$ws.Range('E46').Value = '  +2.58%  '
$ws.Range('E47').Value = '  -4.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.015'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4404'
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.062'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05168'
$ws.Range('E51').Value = '  +0.41%  '
